$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 3552.9524
$ws.Range("I19").Value = 1648.75
$ws.Range("J19").Value = 4724.769
$ws.Range("K19").Value = 1648.75
$ws.Range("L19").Value = 4724.769
$ws.Range("M19").Value = -1473.75
$ws.Range("N19").Value = -5074.769

$ws.Range("H64").Value = 5434.909
$ws.Range("I64").Value = 3723.125
$ws.Range("K64").Value = 3723.125
$ws.Range("M64").Value = -3475.125

$ws.Range("H67").Value = 5434.909
$ws.Range("I67").Value = 3723.125
$ws.Range("K67").Value = 3723.125
$ws.Range("M67").Value = -2865.125

$ws.Range("H88").Value = 8970.5
$ws.Range("J88").Value = 11572.286
$ws.Range("L88").Value = 11572.286
$ws.Range("N88").Value = -12384.286

$ws.Range("H91").Value = 8970.5
$ws.Range("J91").Value = 11572.286
$ws.Range("L91").Value = 11572.286
$ws.Range("N91").Value = -14380.286

$ws.Range("H111").Value = 831.2857
$ws.Range("I111").Value = 857.25
$ws.Range("J111").Value = 796.6667
$ws.Range("K111").Value = 2571.75
$ws.Range("L111").Value = 2390.0001
$ws.Range("M111").Value = 495.25
$ws.Range("N111").Value = -8524.000100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1407.2222
$ws.Range("I2").Value = 1149.4286
$ws.Range("J2").Value = 2309.5
$ws.Range("K2").Value = 1149.4286
$ws.Range("L2").Value = 2309.5
$ws.Range("M2").Value = -1036.4286
$ws.Range("N2").Value = -2535.5

$ws.Range("H74").Value = 408786.6
$ws.Range("I74").Value = 751275
$ws.Range("K74").Value = 751275
$ws.Range("M74").Value = -750401

$ws.Range("H77").Value = 408786.6
$ws.Range("I77").Value = 751275
$ws.Range("K77").Value = 3756375
$ws.Range("M77").Value = -3752007

$ws.Range("H116").Value = 1407.2222
$ws.Range("I116").Value = 1149.4286
$ws.Range("J116").Value = 2309.5
$ws.Range("K116").Value = 1149.4286
$ws.Range("L116").Value = 2309.5
$ws.Range("M116").Value = 1144.5714
$ws.Range("N116").Value = -6897.5

$ws.Range("H132").Value = 872.125
$ws.Range("I132").Value = 769.26
$ws.Range("K132").Value = 2307.78
$ws.Range("M132").Value = 222.2200000000003

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1407.2222
$ws.Range("I3").Value = 1149.4286
$ws.Range("J3").Value = 2309.5
$ws.Range("K3").Value = 1149.4286
$ws.Range("L3").Value = 2309.5
$ws.Range("M3").Value = -1035.4286
$ws.Range("N3").Value = -2537.5

$ws.Range("H86").Value = 1825.5
$ws.Range("I86").Value = 1871.2
$ws.Range("J86").Value = 1749.3334
$ws.Range("K86").Value = 1871.2
$ws.Range("L86").Value = 1749.3334
$ws.Range("M86").Value = -748.2
$ws.Range("N86").Value = -3995.3334

$ws.Range("H89").Value = 1825.5
$ws.Range("I89").Value = 1871.2
$ws.Range("J89").Value = 1749.3334
$ws.Range("K89").Value = 9356
$ws.Range("L89").Value = 8746.666999999999
$ws.Range("M89").Value = -3740
$ws.Range("N89").Value = -19978.667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 874.17645
$ws.Range("I22").Value = 576.1
$ws.Range("J22").Value = 1300
$ws.Range("K22").Value = 576.1
$ws.Range("L22").Value = 1300
$ws.Range("M22").Value = -226.1
$ws.Range("N22").Value = -2000

$ws.Range("H31").Value = 1925088.6
$ws.Range("J31").Value = 4278.8237
$ws.Range("L31").Value = 4278.8237
$ws.Range("N31").Value = -4868.8237

$ws.Range("H34").Value = 1925088.6
$ws.Range("J34").Value = 4278.8237
$ws.Range("L34").Value = 4278.8237
$ws.Range("N34").Value = -4682.8237

$ws.Range("H62").Value = 5888.7
$ws.Range("I62").Value = 6157
$ws.Range("J62").Value = 5486.25
$ws.Range("K62").Value = 6157
$ws.Range("L62").Value = 5486.25
$ws.Range("M62").Value = -5533
$ws.Range("N62").Value = -6734.25

$ws.Range("H65").Value = 5888.7
$ws.Range("I65").Value = 6157
$ws.Range("J65").Value = 5486.25
$ws.Range("K65").Value = 30785
$ws.Range("L65").Value = 27431.25
$ws.Range("M65").Value = -27665
$ws.Range("N65").Value = -33671.25

$ws.Range("H105").Value = 4165.1665
$ws.Range("I105").Value = 4497
$ws.Range("K105").Value = 4497
$ws.Range("M105").Value = -2750

$ws.Range("H134").Value = 1413.2
$ws.Range("I134").Value = 1157
$ws.Range("K134").Value = 3471
$ws.Range("M134").Value = -936

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 4153.846
$ws.Range("I51").Value = 1333.3334
$ws.Range("K51").Value = 4000.0002
$ws.Range("M51").Value = -3540.0002

$ws.Range("H97").Value = 1044.1923
$ws.Range("I97").Value = 690.8333
$ws.Range("K97").Value = 2072.4999
$ws.Range("M97").Value = -1576.4999

$ws.Range("H122").Value = 902.5263
$ws.Range("I122").Value = 388.42856
$ws.Range("J122").Value = 1202.4166
$ws.Range("K122").Value = 3495.85704
$ws.Range("L122").Value = 10821.7494
$ws.Range("M122").Value = -1045.85704
$ws.Range("N122").Value = -15721.7494

$ws.Range("H129").Value = 2984.8
$ws.Range("I129").Value = 1559.6666
$ws.Range("K129").Value = 4678.9998
$ws.Range("M129").Value = 321.0002000000004

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5575.5835
$ws.Range("I122").Value = 5400.5557
$ws.Range("J122").Value = 6100.6665
$ws.Range("K122").Value = 16201.6671
$ws.Range("L122").Value = 18301.9995
$ws.Range("M122").Value = -13751.6671
$ws.Range("N122").Value = -23201.9995

$ws.Range("H132").Value = 4018.8372
$ws.Range("I132").Value = 4069.1562
$ws.Range("J132").Value = 3872.4546
$ws.Range("K132").Value = 12207.4686
$ws.Range("L132").Value = 11617.3638
$ws.Range("M132").Value = -9677.4686
$ws.Range("N132").Value = -16677.3638

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 862.53845
$ws.Range("I16").Value = 873.8182
$ws.Range("J16").Value = 800.5
$ws.Range("K16").Value = 873.8182
$ws.Range("L16").Value = 800.5
$ws.Range("M16").Value = -703.8182
$ws.Range("N16").Value = -1140.5

$ws.Range("H22").Value = 1787.5
$ws.Range("J22").Value = 1757.1428
$ws.Range("L22").Value = 1757.1428
$ws.Range("N22").Value = -2347.1428

$ws.Range("H27").Value = 1787.5
$ws.Range("J27").Value = 1757.1428
$ws.Range("L27").Value = 1757.1428
$ws.Range("N27").Value = -1971.1428

$ws.Range("H31").Value = 1994
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").Value = $null

$ws.Range("H46").Value = 3463.4666
$ws.Range("I46").Value = 1269
$ws.Range("J46").Value = 4560.7
$ws.Range("K46").Value = 1269
$ws.Range("L46").Value = 4560.7
$ws.Range("M46").Value = -1081
$ws.Range("N46").Value = -4936.7

$ws.Range("H132").Value = 3042.6956
$ws.Range("I132").Value = 2554
$ws.Range("J132").Value = 3678
$ws.Range("K132").Value = 7662
$ws.Range("L132").Value = 11034
$ws.Range("M132").Value = -5132
$ws.Range("N132").Value = -16094

$ws.Range("H136").Value = 3578.4546
$ws.Range("I136").Value = 3262.5557
$ws.Range("K136").Value = 9787.667099999999
$ws.Range("M136").Value = -7237.667099999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 8830.5
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").Value = $null

$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("M9").Value = $null

$ws.Range("H126").Value = 253408.9
$ws.Range("I126").Value = 2561.7273
$ws.Range("J126").Value = 559999.9
$ws.Range("K126").Value = 7685.1819
$ws.Range("L126").Value = 1679999.7
$ws.Range("M126").Value = -5215.1819
$ws.Range("N126").Value = -1684939.7

$ws.Range("H136").Value = 15629.881
$ws.Range("I136").Value = 16875.264
$ws.Range("K136").Value = 50625.792
$ws.Range("M136").Value = -48075.792
